$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Write all cell values in the exact order the original author must have
#    typed them (this governs the order new entries are appended to the
#    shared-strings table).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "CRS"
$ws.Range("B2").Value = "SRS"
$ws.Range("C2").Value = "DESIGN"
$ws.Range("D2").Value = "CODE"
$ws.Range("E2").Value = "TEST"
$ws.Range("F2").Value = "RELATED SRS"

$ws.Range("A3").Value = "CRS_00_6"
$ws.Range("B3").Value = "SRS_00_1"
$ws.Range("F3").ClearContents()

$ws.Range("A4").Value = "CRS_00_1"
$ws.Range("B4").Value = "SRS_00_2"
$ws.Range("F4").ClearContents()

$ws.Range("A5").Value = "CRS_00_2, CRS_00_5"
$ws.Range("B5").Value = "SRS_00_3"
$ws.Range("F5").Value = "SRS_00_2"

$ws.Range("A6").Value = "CRS_00_3, CRS_00_4"
$ws.Range("B6").Value = "SRS_00_4"

$ws.Range("B9").Value = "SRS_00_7"

$ws.Range("A7").Value = "CRS_00_7, CRS_00_8, CRS_00_9, CRS_00_11"
$ws.Range("B7").Value = "SRS_00_5"

$ws.Range("A9").Value = "CRS_00_10"

$ws.Range("A8").Value = "CRS_00_12, CRS_00_13, CRS_00_14"
$ws.Range("B8").Value = "SRS_00_6"
$ws.Range("F8").ClearContents()

$ws.Range("F6").Value = "SRS_00_2, SRS_00_3"
$ws.Range("F7").Value = "SRS_00_2, SRS_00_3, SRS_00_4"

$ws.Range("A1").Value = "RTM FOR FAN CONTROLLER SYSTEM"

# ---------------------------------------------------------------------------
# 2) Formatting.
#    Most cells already carry the correct visual formatting (font / fill /
#    border / alignment) from the original file, so only the cells whose
#    *visual* format actually changes need to be touched here.
# ---------------------------------------------------------------------------

# Build the new "grouped CRS" look (14pt font, centered, thin top+right
# border only) once on a scratch cell, then fan it out to A3:A5.
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1").Borders.Item(10).LineStyle = 1
$ws.Range("Z1").Borders.Item(8).LineStyle = 1
$ws.Range("Z1").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# F7 loses the light-blue highlight fill it had before (now plain bordered
# cell with text), matching the look already used by F3:F6/F8.
$ws.Range("F3").Copy()
$ws.Range("F7").PasteSpecial(-4122)

# Brand-new row 9: A9/B9 get the highlighted (light-blue) look used by the
# other grouped rows; C9:F9 get the plain bordered (no text) look.
$ws.Range("A7").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C9:F9").PasteSpecial(-4122)

$ws.Rows.Item(9).RowHeight = 18.75

# ---------------------------------------------------------------------------
# 3) Column F is now much wider to fit the longer "RELATED SRS" lists.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 36.5

# ---------------------------------------------------------------------------
# 4) Selection moves to F14.
# ---------------------------------------------------------------------------
$ws.Range("F14").Select()
